$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Tonquim"
$ws.Range("A15").Value = "Suíça"
$ws.Range("A16").Value = "Japão"
$ws.Range("A18").Value = "Europa"
$ws.Range("A19").Value = "Alsácia"
$ws.Range("A21").Value = "Lituânia"
$ws.Range("A22").Value = "México"
$ws.Range("A23").Value = "Coreia"
$ws.Range("A25").Value = "Dalmácia"
$ws.Range("A26").Value = "Inglaterra"
